# update-24 Dec 2024 -01
#
# Fills in the missing "Dev branch" (column D) values for the last few
# rows of the lookup table on Sheet1, and moves the saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate() | Out-Null

# --- New data for column D (Dev branch) ------------------------------
# Order matters: new shared strings are appended in the order they are
# first written, so D15 (a brand-new string) must be written before
# D13 (also a brand-new string) to land them at shared-string indices
# 61 and 62 respectively.

# Row 15 - stUEmail for Dev branch (new hyperlinked e-mail address)
$ws.Range("D15").Value = "bhupesh+DevStationUser1@atinatechnology.in"
$ws.Hyperlinks.Add($ws.Range("D15"), "mailto:bhupesh+DevStationUser1@atinatechnology.in") | Out-Null
$ws.Range("D15").Style = $ws.Range("C15").Style

# Row 13 - St1AttributeName for Dev branch (plain text, no hyperlink)
$ws.Range("D13").Value = "BH1A1"

# Row 16 - ff1UEmail for Dev branch (reuses the Testing-branch address)
$ws.Range("D16").Value = $ws.Range("C16").Value2
$ws.Hyperlinks.Add($ws.Range("D16"), "mailto:bhupesh+TestingFirefighter1@atinatechnology.in") | Out-Null
$ws.Range("D16").Style = $ws.Range("C15").Style

# Row 17 - ff2UEmail for Dev branch (reuses the Testing-branch address)
$ws.Range("D17").Value = $ws.Range("C17").Value2
$ws.Hyperlinks.Add($ws.Range("D17"), "mailto:bhupesh+TestingFirefighter2@atinatechnology.in") | Out-Null
$ws.Range("D17").Style = $ws.Range("C17").Style

# --- Saved selection ----------------------------------------------------
$ws.Range("C2").Select() | Out-Null
